$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("E7").Value = "Average Among the 3 Sites"
